$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2 "2022-07-20" -> "2022-07-22" (keep as literal text; a leading apostrophe
# stops the date-like string from being auto-converted into a date serial).
$ws.Range("D2").Value = "'2022-07-22"

# E2 "20:02" -> "11:36" (plain text assignment keeps it as text, matching
# the original inline-string cell).
$ws.Range("E2").Value = "11:36"

# F2, G2, H2 numeric updates.
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 99.15573139735331
$ws.Range("H2").Value = 100.4527343703501
